$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 332, shifting existing rows 332.. down by one.
$ws.Range("A332").EntireRow.Insert()

# Populate the newly inserted row 332 with the new weekly data point.
# Most columns mirror the row that used to sit at 332 (now shifted to 333);
# only the date and the K/L/M/P price columns carry new values.
$ws.Cells.Item(332, 1).Value = 9
$ws.Cells.Item(332, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(332, 3).Value = "Metropolitana"
$ws.Cells.Item(332, 4).Value = 44627
$ws.Cells.Item(332, 4).NumberFormat = $ws.Cells.Item(333, 4).NumberFormat
$ws.Cells.Item(332, 5).Value = 13
$ws.Cells.Item(332, 6).Value = 100114014
$ws.Cells.Item(332, 7).Value = "Betarraga"
$ws.Cells.Item(332, 8).Value = "Sin especificar"
$ws.Cells.Item(332, 9).Value = "Primera"
$ws.Cells.Item(332, 10).Value = 4300
$ws.Cells.Item(332, 11).Value = 110
$ws.Cells.Item(332, 12).Value = 120
$ws.Cells.Item(332, 13).Value = 115
$ws.Cells.Item(332, 14).Value = "$/unidad"
$ws.Cells.Item(332, 15).Value = "Región Metropolitana"
$ws.Cells.Item(332, 16).Value = 115
$ws.Cells.Item(332, 17).Value = 1
$ws.Cells.Item(332, 18).Value = "Hortaliza"
